$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Function" rows describing the categories feature (adds shared strings 11-15)
$ws.Range("A4").Value = "Allow users to post messages that show up in their app."
$ws.Range("A5").Value = "Allow users to browse existing categories."
$ws.Range("A6").Value = "Allow users to create a new category, if the category they're looking for doesn't already exist."
$ws.Range("A7").Value = "Allow users to select a message to read more."
$ws.Range("A8").Value = "Allow users to comment on/reply to messages."

# New "Status" column (F) for the existing task rows
$ws.Range("F1").Value = "Status"
$ws.Range("F3").Value = "On going"
$ws.Range("F2").Value = "Done"

# Widen column A to fit the new, longer text (~77.25 chars, like the authored file)
$ws.Columns.Item(1).ColumnWidth = 76.5

# Page setup (portrait, paper size 9 = A4)
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

# Move the active selection like in the authored workbook
$ws.Range("C12").Select() | Out-Null
